# "allow empty cells when evaluate"
# Root-cause data edit: A4 goes from 4 -> 20, which ripples through every
# formula that (directly or via the B/C shared-formula chains and LINEST)
# depends on column A. D9's OFFSET formula is widened (A4*2 instead of A4)
# so its SUM range reaches past the populated data into empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 20
$ws.Range("D9").Formula = "=SUM(A1:OFFSET(A2,A4*2,0))"

# Reflect the author's final selection (cell below the data, A19) so the
# saved sheetView carries the same <selection> as the target workbook.
$ws.Range("A19").Select()
